# Automatische test-sync: 2025-06-30 19:51:50
# Adds Testmail #7 to the "Logs" sheet and updates the "Dashboard" sheet's
# category counts / ordering to reflect the new entry.

$wb = $excel.ActiveWorkbook
$logs = $wb.Worksheets.Item("Logs")
$dashboard = $wb.Worksheets.Item("Dashboard")

# --- Logs sheet: append row 8 with the new test mail data ---------------
$logs.Range("A8").Value = "Ik ben niet tevreden met mijn bestelling."
$logs.Range("B8").Value = "mailmind.test@zohomail.eu"
$logs.Range("C8").Value = "Testmail #7: Ik ben niet tevreden met mijn bestelling."
$logs.Range("D8").Value = "Retour / Terugbetaling"
$logs.Range("E8").Value = "Beste klant,`nBedankt voor uw bericht. Wat vervelend om te horen dat u niet tevreden bent met uw bestelling. Om u beter van dienst te kunnen zijn, zou ik graag meer details willen weten over wat er precies niet naar wens is gegaan. Kunt u mogelijk informatie geven over het specifieke product of de reden waarom u niet tevreden bent? Op die manier kunnen we het probleem verder onderzoeken en een passende oplossing bieden.`nAlvast bedankt voor uw medewerking.`nMet vriendelijke groet,`n[Naam]  `nE-mailassistent  `n[Bedrijfsnaam]"
$logs.Range("F8").Value = "2025-06-30 19:51:49"
$logs.Range("G8").Value = "Ja"
$logs.Range("H8").Value = "Nee"
$logs.Range("I8").Value = "Ja"
$logs.Range("J8").Value = "Nee"

# Re-fit the row height so it matches the sheet's default (avoids a stray
# customHeight flag being written for the newly populated row).
$logs.Rows.Item(8).AutoFit()

# --- Logs sheet: extend the conditional formatting ranges to include row 8
$ranges = @("D2:D7", "G2:G7", "H2:H7", "I2:I7", "J2:J7")
foreach ($rng in $ranges) {
    $col = $rng.Substring(0, 1)
    $newRange = "$col" + "2:" + "$col" + "8"
    $fcs = $logs.Range($rng).FormatConditions
    for ($i = 1; $i -le $fcs.Count; $i++) {
        $fcs.Item($i).ModifyAppliesToRange($logs.Range($newRange))
    }
}

# --- Dashboard sheet: category order swap + updated counts --------------
$dashboard.Range("A2").Value = "Retour / Terugbetaling"
$dashboard.Range("B2").Value = 3
$dashboard.Range("A3").Value = "Productinformatie"
$dashboard.Range("B3").Value = 2
